$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 86
$prev = 85

# Copy the formatting from the previous (last existing) data row so the new
# row's styles (bordered/bold index column, date-formatted match-date column)
# match the rest of the table exactly.
$ws.Range("A$prev").Copy()
$ws.Range("A$row").PasteSpecial(-4122)
$ws.Range("E$prev").Copy()
$ws.Range("E$row").PasteSpecial(-4122)

$ws.Cells.Item($row, 1).Value = 85
$ws.Cells.Item($row, 2).Value = "azerbaijan"
$ws.Cells.Item($row, 3).Value = "premier-league"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45282.54166666666
$ws.Cells.Item($row, 6).Value = "Gabala"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = "Sabah Baku"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 3.09
$ws.Cells.Item($row, 11).Value = "21/12/2023 01:12"
$ws.Cells.Item($row, 12).Value = 3.94
$ws.Cells.Item($row, 13).Value = "22/12/2023 12:57"
$ws.Cells.Item($row, 14).Value = 3.51
$ws.Cells.Item($row, 15).Value = "21/12/2023 01:12"
$ws.Cells.Item($row, 16).Value = 4.17
$ws.Cells.Item($row, 17).Value = "22/12/2023 12:59"
$ws.Cells.Item($row, 18).Value = 2.01
$ws.Cells.Item($row, 19).Value = "21/12/2023 01:12"
$ws.Cells.Item($row, 20).Value = 1.73
$ws.Cells.Item($row, 21).Value = "22/12/2023 12:59"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/gabala-sabah-baku/WU5AgyAQ/"
